$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for new rows 101-109 (TimeStamp, then 14 numeric columns B..O)
$rows = @(
    @(45726.720196759263, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4702),
    @(45726.721956018519, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4702),
    @(45727.226168981484, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725),
    @(45727.233865740738, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725),
    @(45727.237384259257, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725),
    @(45727.24659722222,  10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725),
    @(45727.298449074071, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725),
    @(45727.304571759261, 10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725),
    @(45727.33902777778,  10, 6, 240, 426, 402, 476, 3432, 476, 2026, 208, 417, 30, 3683, 4725)
)

$startRow = 101
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
    # Copy style formatting from the row above so number formats / alignment match existing data
    $ws.Range($ws.Cells.Item($r - 1, 1), $ws.Cells.Item($r - 1, 15)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 15)).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
